$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp in the title cell
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 13:20"

# Re-sorted / refreshed country statistics (Excel COM is 1-based row/col)
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 215357
$ws.Cells.Item(4, 3).Value = 354
$ws.Cells.Item(4, 4).Value = 8878
$ws.Cells.Item(4, 5).Value = 201366
$ws.Cells.Item(4, 6).Value = 5005
$ws.Cells.Item(4, 7).Value = 11
$ws.Cells.Item(4, 8).Value = 5113

$ws.Cells.Item(16, 1).Value = "Austria"
$ws.Cells.Item(16, 2).Value = 10892
$ws.Cells.Item(16, 3).Value = 181
$ws.Cells.Item(16, 4).Value = 1749
$ws.Cells.Item(16, 5).Value = 8985
$ws.Cells.Item(16, 6).Value = 227
$ws.Cells.Item(16, 7).Value = 12
$ws.Cells.Item(16, 8).Value = 158

$ws.Cells.Item(28, 1).Value = "Dinamarca"
$ws.Cells.Item(28, 2).Value = 3355
$ws.Cells.Item(28, 3).Value = 248
$ws.Cells.Item(28, 4).Value = 1089
$ws.Cells.Item(28, 5).Value = 2143
$ws.Cells.Item(28, 6).Value = 153
$ws.Cells.Item(28, 7).Value = 19
$ws.Cells.Item(28, 8).Value = 123

$ws.Cells.Item(75, 1).Value = "Principado de Andorra"
$ws.Cells.Item(75, 2).Value = 428
$ws.Cells.Item(75, 3).Value = 38
$ws.Cells.Item(75, 4).Value = 10
$ws.Cells.Item(75, 5).Value = 403
$ws.Cells.Item(75, 6).Value = 12
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 15

$ws.Cells.Item(76, 1).Value = "Eslovaquia"
$ws.Cells.Item(76, 2).Value = 426
$ws.Cells.Item(76, 3).Value = 26
$ws.Cells.Item(76, 4).Value = 3
$ws.Cells.Item(76, 5).Value = 422
$ws.Cells.Item(76, 6).Value = 3
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 1

$ws.Cells.Item(77, 1).Value = "Tunez"
$ws.Cells.Item(77, 2).Value = 423
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 5
$ws.Cells.Item(77, 5).Value = 406
$ws.Cells.Item(77, 6).Value = 10
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 12

$ws.Cells.Item(78, 1).Value = "Moldavia"
$ws.Cells.Item(78, 2).Value = 423
$ws.Cells.Item(78, 3).Value = 0
$ws.Cells.Item(78, 4).Value = 23
$ws.Cells.Item(78, 5).Value = 395
$ws.Cells.Item(78, 6).Value = 65
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 5

$ws.Cells.Item(79, 1).Value = "Kazajistan"
$ws.Cells.Item(79, 2).Value = 402
$ws.Cells.Item(79, 3).Value = 22
$ws.Cells.Item(79, 4).Value = 27
$ws.Cells.Item(79, 5).Value = 372
$ws.Cells.Item(79, 6).Value = 6
$ws.Cells.Item(79, 7).Value = 0
$ws.Cells.Item(79, 8).Value = 3

$ws.Cells.Item(80, 1).Value = "Azerbaiyan"
$ws.Cells.Item(80, 2).Value = 400
$ws.Cells.Item(80, 3).Value = 41
$ws.Cells.Item(80, 4).Value = 26
$ws.Cells.Item(80, 5).Value = 369
$ws.Cells.Item(80, 6).Value = 7
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 5

$ws.Cells.Item(95, 1).Value = "Vietnam"
$ws.Cells.Item(95, 2).Value = 227
$ws.Cells.Item(95, 3).Value = 9
$ws.Cells.Item(95, 4).Value = 75
$ws.Cells.Item(95, 5).Value = 152
$ws.Cells.Item(95, 6).Value = 3
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 0

$ws.Cells.Item(139, 1).Value = "Zambia"
$ws.Cells.Item(139, 2).Value = 39
$ws.Cells.Item(139, 3).Value = 3
$ws.Cells.Item(139, 4).Value = 0
$ws.Cells.Item(139, 5).Value = 38
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = 1

$ws.Cells.Item(140, 1).Value = "Puerto Rico"
$ws.Cells.Item(140, 2).Value = 39
$ws.Cells.Item(140, 3).Value = 0
$ws.Cells.Item(140, 4).Value = 1
$ws.Cells.Item(140, 5).Value = 36
$ws.Cells.Item(140, 6).Value = 0
$ws.Cells.Item(140, 7).Value = 0
$ws.Cells.Item(140, 8).Value = 2

$ws.Cells.Item(141, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(141, 2).Value = 37
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(141, 4).Value = 0
$ws.Cells.Item(141, 5).Value = 37
$ws.Cells.Item(141, 6).Value = 1
$ws.Cells.Item(141, 7).Value = 0
$ws.Cells.Item(141, 8).Value = 0

$ws.Cells.Item(191, 1).Value = "Cabo Verde"
$ws.Cells.Item(191, 2).Value = 6
$ws.Cells.Item(191, 3).Value = 0
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 5
$ws.Cells.Item(191, 6).Value = 0
$ws.Cells.Item(191, 7).Value = 0
$ws.Cells.Item(191, 8).Value = 1

$ws.Cells.Item(192, 1).Value = "Nepal"
$ws.Cells.Item(192, 2).Value = 6
$ws.Cells.Item(192, 3).Value = 1
$ws.Cells.Item(192, 4).Value = 1
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 0
$ws.Cells.Item(192, 7).Value = 0
$ws.Cells.Item(192, 8).Value = 0

$ws.Cells.Item(193, 1).Value = "San Bartolome"
$ws.Cells.Item(193, 2).Value = 6
$ws.Cells.Item(193, 3).Value = 0
$ws.Cells.Item(193, 4).Value = 1
$ws.Cells.Item(193, 5).Value = 5
$ws.Cells.Item(193, 6).Value = 0
$ws.Cells.Item(193, 7).Value = 0
$ws.Cells.Item(193, 8).Value = 0

$ws.Cells.Item(194, 1).Value = "Mauritania"
$ws.Cells.Item(194, 2).Value = 6
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 2
$ws.Cells.Item(194, 5).Value = 3
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 1

$ws.Cells.Item(195, 1).Value = "Nicaragua"
$ws.Cells.Item(195, 2).Value = 5
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 0
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 1

$ws.Cells.Item(197, 1).Value = "Somalia"
$ws.Cells.Item(197, 2).Value = 5
$ws.Cells.Item(197, 3).Value = 0
$ws.Cells.Item(197, 4).Value = 1
$ws.Cells.Item(197, 5).Value = 4
$ws.Cells.Item(197, 6).Value = 0
$ws.Cells.Item(197, 7).Value = 0
$ws.Cells.Item(197, 8).Value = 0
